# Add two new columns (C: NewBuildingName / D: NewFloorname) of sample
# data to the "manageBuilding" sheet, matching the header formatting of
# the existing columns, and leave the selection on B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manageBuilding")

# New header row values (column C then column D so the shared-string
# table is built up in the same order the sample data appears).
$ws.Range("C1").Value = "NewBuildingName"
$ws.Range("C2").Value = "building6"
$ws.Range("C3").Value = "building10"

$ws.Range("D1").Value = "NewFloorname"
$ws.Range("D2").Value = "floorNo12"
$ws.Range("D3").Value = "FloorNo5"

# Match the yellow header-row fill used by the existing A1/B1 headers.
$ws.Range("C1").Interior.Color = $ws.Range("A1").Interior.Color
$ws.Range("D1").Interior.Color = $ws.Range("B1").Interior.Color

# Give the new columns sensible custom widths (close to the sizing the
# existing columns use for header text of this length).
$ws.Columns.Item(3).ColumnWidth = 16.25
$ws.Columns.Item(4).ColumnWidth = 17

# Leave the active selection on B1.
$ws.Range("B1").Select() | Out-Null
